$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A1").Font.Size = 11
$ws1.Range("A1").Font.Color = 16777215
